$wb = $excel.ActiveWorkbook

# --- Update the daily conversion note on "Hoja1" ---
$ws1 = $wb.Worksheets.Item("Hoja1")

$ws1.Range("A1").Value = "Conversión del día 💰`n✅ Dólar paralelo: 68`n`nBinance`n✅ 1000 Bs = 9.48 = 40868.58 pesos`n✅ 40868.58 pesos = 9.45 = 954.86 Bs`n`nPromedio competencia`n✅ Tasa pesos: 20`n✅ Tasa Bs: 20`n✅ % Ganancia: 20%"

# --- Update the rate figures on "tasas" ---
$ws2 = $wb.Worksheets.Item("tasas")

$ws2.Range("N10").Value = 105.46
$ws2.Range("O10").Value = 4310
$ws2.Range("N12").Value = 4325
$ws2.Range("O12").Value = 101.05
